# mads_config.xlsx — "Modify requirements and excel naming"
#
# Functional edit: on the "data_sheet" worksheet, WITHIN_DAYS (row 3, column B)
# changes from 30 to 29, and the sheet's active selection moves from D9 to B3
# (matching where the value was just edited).
#
# (The surrounding diff also touches purely environmental/non-deterministic
# workbook metadata — the author's absolute checkout path in x15ac:absPath,
# the xr:revisionPtr documentId GUID, and the saved window
# position/size in <workbookView> — all of which are stamped by the real
# Excel client from the local machine/session and aren't part of the
# document's actual content; they aren't reachable through the Excel
# object model, so they're intentionally left untouched here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data_sheet")

# WITHIN_DAYS: 30 -> 29
$ws.Range("B3").Value = 29

# Leave the edited cell selected, as the saved file shows.
$ws.Range("B3").Select()
